# Add a new bullet to the Freddie Mac experience block: a new
# "Constructed an optimization solver ..." line right after the
# "... visualizations to support analysis and performance evaluation"
# bullet and before the "Collaborated with team members ..." bullet.

$d = $word.ActiveDocument

$anchorText = "visualizations to support analysis and performance evaluation"
$newText = "Constructed an optimization solver to reduce time to determine loan weights from over 2 hours to under 2 minutes"

# Locate the paragraph that ends with the anchor text.
$anchorPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like ("*" + $anchorText + "*")) {
        $anchorPara = $p
    }
}

# Insert a brand-new paragraph right after it; it inherits the anchor
# paragraph's list/style/run formatting, matching the surrounding bullets.
$anchorPara.Range.InsertParagraphAfter()

$anchorIndex = $anchorPara.Index
$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newPara.Range.Text = $newText
